# Rename model "Region" to "Zone" for clarity, and make the renamed sheet
# the active one (mirroring the author's last-saved selection state).

$wb = $excel.ActiveWorkbook

# 1. Rename the "Region" worksheet to "Zone".
$zone = $wb.Worksheets.Item("Region")
$zone.Name = "Zone"

# 2. The previously-active sheet ("Line") loses the tab-selected flag while
#    keeping its existing selection; the renamed "Zone" sheet becomes active
#    with a new selection.
$zone.Activate() | Out-Null
$zone.Range("D25").Select() | Out-Null
